$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 332 (the post "「楽観主義」..." entry) - all subsequent rows shift up by one.
$ws.Rows.Item(332).Delete()
